# Update cryptos list data (values + two row swaps) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to stay text even though the new value parses as a number
# (Excel auto-converts numeric-looking strings unless the cell is pre-formatted as text)
$textCells = @("D5", "D6", "D9", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D31", "D36", "D37", "D42", "D43", "D45", "D46", "D47", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the cell values
$ws.Range("D2").Value = "72.316.94"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "2.641.33"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "602.75"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "180.77"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "0.177"
$ws.Range("E9").Value = "  +6.90%  "
$ws.Range("D10").Value = "2.640.49"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "3.124.91"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "72.168.84"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "26.62"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "2.641.64"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").Value = "11.97"
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").Value = "7.93"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "378.53"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "2.07"
$ws.Range("E23").Value = "  +10.40%  "
$ws.Range("D24").Value = "73.08"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "4.40"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("D28").Value = "2.777.83"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "0.0₃0959"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "524.47"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "164.02"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "19.35"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("E41").Value = "  -3.99%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "5.10"
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "39.35"
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").Value = "151.66"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  -4.35%  "

# Restore default (Normal) style on the cells we temporarily formatted as text
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
